$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 59, shifting existing rows 59:65 down to 60:66
$ws.Rows.Item(59).Insert()

# Populate the newly inserted row 59 with the new record's data
$ws.Range("A59").Value = 6
$ws.Range("B59").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C59").Value = "Metropolitana"
$ws.Range("D59").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D59").Value = 45275
$ws.Range("E59").Value = 13
$ws.Range("F59").Value = "Fruta"
$ws.Range("G59").Value = 100102
$ws.Range("H59").Value = "Cítricos"
$ws.Range("I59").Value = 100102006
$ws.Range("J59").Value = "Pomelo"
$ws.Range("K59").Value = "Start Ruby"
$ws.Range("L59").Value = "Primera"
$ws.Range("M59").Value = 20
$ws.Range("N59").Value = 190000
$ws.Range("O59").Value = 190000
$ws.Range("P59").Value = 190000
$ws.Range("Q59").Value = '$/bins (350 kilos)'
$ws.Range("R59").Value = "Provincia de Limarí"
$ws.Range("S59").Value = 543
$ws.Range("T59").Value = 350

Write-Host "Done"
